# Append a new tag ("ถูกสาม") with its keyword variants as new rows at the
# bottom of the training data sheet (mirrors the existing tag/keyword
# layout already present in columns A/B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tag = "ถูกสาม"
$keywords = @(
    "ปีโป้",
    "ปีโป๊",
    "ปีโป้หลากสี",
    "ปีปีโป้ปะปะปีปีโป้",
    "ปีโป้ไง",
    "ปีโป้ครับ"
)

# The sheet already uses a blank separator row between each tag group
# (e.g. row 269 between the "ถูกหนึ่ง" and "ถูกสอง" groups). Keep that same
# style: leave one blank row after the last existing entry (row 274) and
# start the new "ถูกสาม" group at row 276.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row = $lastRow + 2

foreach ($kw in $keywords) {
    $ws.Cells.Item($row, 1).Value = $tag
    $ws.Cells.Item($row, 2).Value = $kw
    $row = $row + 1
}

# Match the author's final selection/scroll position recorded in the sheet.
$ws.Activate()
$lastCell = $ws.Cells.Item($row - 1, 2)
$lastCell.Select()
